$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3488604574705505
$ws.Range("C2").Value = 0.06922633686004076
$ws.Range("D2").Value = 0.0786235355372753
$ws.Range("E2").Value = 0.4135961380017221
$ws.Range("G2").Value = 0.002437001023303067
$ws.Range("I2").Value = 0.6179804247686214
$ws.Range("K2").Value = 0.3928098750408537
$ws.Range("O2").Value = 2.939250834118553
$ws.Range("B3").Value = 0.3089029672826484
$ws.Range("C3").Value = 0.06034598482280273
$ws.Range("D3").Value = 0.07130418711329867
$ws.Range("E3").Value = 0.3608725185316786
$ws.Range("G3").Value = 0.002439999921773311
$ws.Range("I3").Value = 0.6215252185763802
$ws.Range("K3").Value = 0.3455162402882763
$ws.Range("O3").Value = 2.938201290481601
$ws.Range("B4").Value = 0.2843885206028176
$ws.Range("C4").Value = 0.05487995987675731
$ws.Range("D4").Value = 0.06684519773406805
$ws.Range("E4").Value = 0.3285912352611149
$ws.Range("G4").Value = 0.002441937328537203
$ws.Range("I4").Value = 0.6240811998009725
$ws.Range("K4").Value = 0.316476674743285
$ws.Range("O4").Value = 2.939490406596235
$ws.Range("B5").Value = 0.2744039641116274
$ws.Range("C5").Value = 0.05264911350766965
$ws.Range("D5").Value = 0.06503693181619496
$ws.Range("E5").Value = 0.3154575345480595
$ws.Range("G5").Value = 0.002442751073565055
$ws.Range("I5").Value = 0.625217999181892
$ws.Range("K5").Value = 0.3046428383502189
$ws.Range("O5").Value = 2.940500720834279
$ws.Range("B6").Value = 0.2727463653448012
$ws.Range("C6").Value = 0.05227847763528359
$ws.Range("D6").Value = 0.06473720199592492
$ws.Range("E6").Value = 0.31327792190244
$ws.Range("G6").Value = 0.00244288766136757
$ws.Range("I6").Value = 0.6254125087569378
$ws.Range("K6").Value = 0.302677849017698
$ws.Range("O6").Value = 2.940697739429027
$ws.Range("B7").Value = 0.2842538435409097
$ws.Range("C7").Value = 0.0548498876736403
$ws.Range("D7").Value = 0.06682077519144514
$ws.Range("E7").Value = 0.328414026107879
$ws.Range("G7").Value = 0.0024419482048059
$ws.Range("I7").Value = 0.6240961458205661
$ws.Range("K7").Value = 0.316317079000612
$ws.Range("O7").Value = 2.939502069974026
$ws.Range("B8").Value = 0.3350792386076762
$ws.Range("C8").Value = 0.06616717222738089
$ws.Range("D8").Value = 0.07609251889113011
$ws.Range("E8").Value = 0.3953970240016247
$ws.Range("G8").Value = 0.002438015151102247
$ws.Range("I8").Value = 0.6191237891052133
$ws.Range("K8").Value = 0.3765034149490134
$ws.Range("O8").Value = 2.938486903741392
$ws.Range("B9").Value = 0.4348932631119453
$ws.Range("C9").Value = 0.08825535744827562
$ws.Range("D9").Value = 0.09455506540996339
$ws.Range("E9").Value = 0.5275598072005607
$ws.Range("G9").Value = 0.002431061147318442
$ws.Range("I9").Value = 0.6123934197113527
$ws.Range("K9").Value = 0.4945132579337894
$ws.Range("O9").Value = 2.951897036977698
$ws.Range("B10").Value = 0.5083095854439819
$ws.Range("C10").Value = 0.1044236370520935
$ws.Range("D10").Value = 0.1082948454944557
$ws.Range("E10").Value = 0.625280432431822
$ws.Range("G10").Value = 0.002426409541789945
$ws.Range("I10").Value = 0.6093036924781856
$ws.Range("K10").Value = 0.5812056910332899
$ws.Range("O10").Value = 2.971226239622069
$ws.Range("B11").Value = 0.5417259162998675
$ws.Range("C11").Value = 0.1117669938820427
$ws.Range("D11").Value = 0.1145845000201717
$ws.Range("E11").Value = 0.6698989831981379
$ws.Range("G11").Value = 0.002424391676371753
$ws.Range("I11").Value = 0.6083039261767382
$ws.Range("K11").Value = 0.6206429416421599
$ws.Range("O11").Value = 2.982096840127156
$ws.Range("B12").Value = 0.5543823056831343
$ws.Range("C12").Value = 0.1145461001639774
$ws.Range("D12").Value = 0.1169719371382598
$ws.Range("E12").Value = 0.6868207006622669
$ws.Range("G12").Value = 0.002423641599095916
$ws.Range("I12").Value = 0.6079839296192517
$ws.Range("K12").Value = 0.6355767243837818
$ws.Range("O12").Value = 2.986513497538596
$ws.Range("B13").Value = 0.5516564269991022
$ws.Range("C13").Value = 0.1139476436316613
$ws.Range("D13").Value = 0.1164575068084019
$ws.Range("E13").Value = 0.6831751348615853
$ws.Range("G13").Value = 0.002423802518214553
$ws.Range("I13").Value = 0.6080502367444609
$ws.Range("K13").Value = 0.632360482969375
$ws.Range("O13").Value = 2.985548917082497
$ws.Range("B14").Value = 0.5427671200397413
$ws.Range("C14").Value = 0.1119956658761225
$ws.Range("D14").Value = 0.1147808019193661
$ws.Range("E14").Value = 0.6712906181240044
$ws.Range("G14").Value = 0.002424329685935519
$ws.Range("I14").Value = 0.6082764241665259
$ws.Range("K14").Value = 0.621871559660292
$ws.Range("O14").Value = 2.982454176558178
$ws.Range("B15").Value = 0.5373224581851161
$ws.Range("C15").Value = 0.1107998069692826
$ws.Range("D15").Value = 0.113754512297362
$ws.Range("E15").Value = 0.6640144049970047
$ws.Range("G15").Value = 0.002424654418476199
$ws.Range("I15").Value = 0.6084226083422664
$ws.Range("K15").Value = 0.6154467498734562
$ws.Range("O15").Value = 2.980597695809251
$ws.Range("B16").Value = 0.5061260982206477
$ws.Range("C16").Value = 0.1039434977587348
$ws.Range("D16").Value = 0.1078845956373158
$ws.Range("E16").Value = 0.6223679881735507
$ws.Range("G16").Value = 0.002426543383477364
$ws.Range("I16").Value = 0.6093772154873065
$ws.Range("K16").Value = 0.5786283559468757
$ws.Range("O16").Value = 2.970557750281444
$ws.Range("B17").Value = 0.4869927157791096
$ws.Range("C17").Value = 0.0997343943533906
$ws.Range("D17").Value = 0.104293688152822
$ws.Range("E17").Value = 0.5968628683563679
$ws.Range("G17").Value = 0.002427727294962482
$ws.Range("I17").Value = 0.6100669361268132
$ws.Range("K17").Value = 0.556041351011487
$ws.Range("O17").Value = 2.964931722162703
$ws.Range("B18").Value = 0.4759894985475057
$ws.Range("C18").Value = 0.09731232953083691
$ws.Range("D18").Value = 0.1022319965652088
$ws.Range("E18").Value = 0.5822084077339156
$ws.Range("G18").Value = 0.002428417494452665
$ws.Range("I18").Value = 0.6105018231135517
$ws.Range("K18").Value = 0.5430499615679594
$ws.Range("O18").Value = 2.961891251513521
$ws.Range("B19").Value = 0.4722643183071398
$ws.Range("C19").Value = 0.09649207065248788
$ws.Range("D19").Value = 0.1015345786282325
$ws.Range("E19").Value = 0.5772492475576598
$ws.Range("G19").Value = 0.002428652774132136
$ws.Range("I19").Value = 0.6106556180420384
$ws.Range("K19").Value = 0.5386513198956209
$ws.Range("O19").Value = 2.960895329364433
$ws.Range("B20").Value = 0.4890293135451032
$ws.Range("C20").Value = 0.1001825743242364
$ws.Range("D20").Value = 0.1046755629396898
$ws.Range("E20").Value = 0.5995763232201767
$ws.Range("G20").Value = 0.002427600309287348
$ws.Range("I20").Value = 0.6099895611440402
$ws.Range("K20").Value = 0.5584457724687297
$ws.Range("O20").Value = 2.965510381470779
$ws.Range("B21").Value = 0.5453780640695243
$ws.Range("C21").Value = 0.112569054080069
$ws.Range("D21").Value = 0.1152731364234825
$ws.Range("E21").Value = 0.6747806774853586
$ws.Range("G21").Value = 0.002424174463340687
$ws.Range("I21").Value = 0.6082083952784885
$ws.Range("K21").Value = 0.6249524198627512
$ws.Range("O21").Value = 2.983355018143556
$ws.Range("B22").Value = 0.5822187629965185
$ws.Range("C22").Value = 0.1206546461209541
$ws.Range("D22").Value = 0.1222323893024537
$ws.Range("E22").Value = 0.7240814296563229
$ws.Range("G22").Value = 0.002422017304169488
$ws.Range("I22").Value = 0.6073859410726996
$ws.Range("K22").Value = 0.6684167729765136
$ws.Range("O22").Value = 2.996767968340038
$ws.Range("B23").Value = 0.5625550693566765
$ws.Range("C23").Value = 0.1163400945734736
$ws.Range("D23").Value = 0.1185150668646315
$ws.Range("E23").Value = 0.6977542924298064
$ws.Range("G23").Value = 0.002423161157202028
$ws.Range("I23").Value = 0.6077935590927623
$ws.Range("K23").Value = 0.6452192754348687
$ws.Range("O23").Value = 2.989448582474637
$ws.Range("B24").Value = 0.4881085775774636
$ws.Range("C24").Value = 0.09997995899948364
$ws.Range("D24").Value = 0.1045029087275395
$ws.Range("E24").Value = 0.598349543135086
$ws.Range("G24").Value = 0.002427657689885588
$ws.Range("I24").Value = 0.610024422904047
$ws.Range("K24").Value = 0.5573587516356611
$ws.Range("O24").Value = 2.965248165348441
$ws.Range("B25").Value = 0.4078759589691856
$ws.Range("C25").Value = 0.08229073925042485
$ws.Range("D25").Value = 0.08952997106145233
$ws.Range("E25").Value = 0.4917065409114088
$ws.Range("G25").Value = 0.00243286168645336
$ws.Range("I25").Value = 0.6138893901710318
$ws.Range("K25").Value = 0.4625901418464196
$ws.Range("O25").Value = 2.946611403764678
